$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11  = -21.678
    12  = -21.544
    15  = -22.062
    27  = -21.776
    28  = -21.889
    31  = -21.849
    32  = -21.742
    36  = -20.365
    38  = -19.829
    46  = -21.873
    54  = -22.002
    55  = -22.184
    56  = -21.933
    67  = -21.577
    69  = -21.503
    72  = -21.689
    73  = -20.137
    83  = -21.988
    86  = -22.135
    91  = -20.887
    93  = -21.452
    99  = -22.118
    104 = -21.335
    105 = -20.241
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
